# Add an "Address" column (new column F) between "Name" (E) and "District"
# (old F, becomes new G). The address text is derived from the second line
# of each person's column-B entry: drop the trailing district segment
# (the text after the last ", ") and concatenate the remaining
# comma-separated segments together (no separators).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-AddressFromBio($bioText) {
    if ($null -eq $bioText) {
        return ""
    }
    $lines = $bioText -split "`n"
    # The normal case is a two line entry: "NAME" then
    # "School, Place, Taluk, District." When that assumption doesn't hold
    # (e.g. an extra embedded line break) no address could be derived.
    if ($lines.Count -ne 2) {
        return ""
    }
    $addrLine = $lines[1]
    $parts = $addrLine -split ", "
    if ($parts.Count -le 1) {
        return ""
    }
    $addr = ($parts[0..($parts.Count - 2)] -join "")
    return $addr
}

# Find the last used row in column A (the SL. NO. column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Insert a new blank column at F; the existing "District" column (F) shifts
# to G, and everything else to the right of the insertion point shifts too.
$ws.Columns.Item(6).Insert()

# Header row.
$ws.Range("F2").Value = "Address"

# Fill in the address for every data row (row 3 through the last row),
# derived from the corresponding column B biography text.
for ($r = 3; $r -le $lastRow; $r++) {
    $bVal = $ws.Range("B" + $r).Value()
    $addr = Get-AddressFromBio $bVal
    if ($addr -ne "") {
        $ws.Range("F" + $r).Value = $addr
    }
}
